$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 11 mirrors the formatting of row 10 (same borders/number formats),
# then gets its own values: a new UASG code, a new org name ("TESTE"),
# and the same acronym as row 10 ("HNBRA").
$ws.Range("A10:C10").Copy()
$ws.Range("A11:C11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("A11").Value = 787701
$ws.Range("B11").Value = "TESTE"
$ws.Range("C11").Value = "HNBRA"
